$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D to hold the "Status" field
$ws.Columns.Item(4).Insert()

# Header row
$ws.Range("D1").Value = "Status"
$ws.Range("E1").Value = "Jan_2026"
$ws.Range("F1").Value = "Dec_2025"
$ws.Range("G1").Value = "Oct_2025"
$ws.Range("H1").Value = "MoM"
$ws.Range("I1").Value = "QoQ"

# Row 2: INE018A01030 - Larsen & Toubro Limited
$ws.Cells.Item(2, 1).Value = "INE018A01030"
$ws.Cells.Item(2, 2).Value = "Larsen & Toubro Limited"
$ws.Cells.Item(2, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(2, 4).Value = "Adding Consistently"
$ws.Cells.Item(2, 5).Value = 9.783791
$ws.Cells.Item(2, 6).Value = 9.613726
$ws.Cells.Item(2, 7).Value = 9.458393
$ws.Cells.Item(2, 8).Value = 0.170065000000001
$ws.Cells.Item(2, 9).Value = 0.3253980000000016

# Row 3: INE002A01018 - Reliance Industries Limited
$ws.Cells.Item(3, 1).Value = "INE002A01018"
$ws.Cells.Item(3, 2).Value = "Reliance Industries Limited"
$ws.Cells.Item(3, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(3, 4).Value = "Reducing Consistently"
$ws.Cells.Item(3, 5).Value = 9.299394
$ws.Cells.Item(3, 6).Value = 9.902957
$ws.Cells.Item(3, 7).Value = 9.342141
$ws.Cells.Item(3, 8).Value = -0.6035630000000012
$ws.Cells.Item(3, 9).Value = -0.04274700000000031

# Row 4: INE775A01035 - Samvardhana Motherson International Ltd
$ws.Cells.Item(4, 1).Value = "INE775A01035"
$ws.Cells.Item(4, 2).Value = "Samvardhana Motherson International Ltd"
$ws.Cells.Item(4, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(4, 4).Value = "Reducing"
$ws.Cells.Item(4, 5).Value = 7.233256
$ws.Cells.Item(4, 6).Value = 7.271142
$ws.Cells.Item(4, 7).Value = 6.367869
$ws.Cells.Item(4, 8).Value = -0.03788600000000031
$ws.Cells.Item(4, 9).Value = 0.8653870000000001

# Row 5: INE814H01029 - Adani Power Limited
$ws.Cells.Item(5, 1).Value = "INE814H01029"
$ws.Cells.Item(5, 2).Value = "Adani Power Limited"
$ws.Cells.Item(5, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(5, 4).Value = "Adding"
$ws.Cells.Item(5, 5).Value = 7.203859
$ws.Cells.Item(5, 6).Value = 7.19971
$ws.Cells.Item(5, 7).Value = 7.921549
$ws.Cells.Item(5, 8).Value = 0.004148999999999958
$ws.Cells.Item(5, 9).Value = -0.7176900000000002

# Row 6: INE406A01037 - Aurobindo Pharma Limited
$ws.Cells.Item(6, 1).Value = "INE406A01037"
$ws.Cells.Item(6, 2).Value = "Aurobindo Pharma Limited"
$ws.Cells.Item(6, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(6, 4).Value = "Adding Consistently"
$ws.Cells.Item(6, 5).Value = 6.905352
$ws.Cells.Item(6, 6).Value = 4.935309
$ws.Cells.Item(6, 7).Value = 4.73556
$ws.Cells.Item(6, 8).Value = 1.970043
$ws.Cells.Item(6, 9).Value = 2.169791999999999

# Row 7: INE758E01017 - Jio Financial Services Limited
$ws.Cells.Item(7, 1).Value = "INE758E01017"
$ws.Cells.Item(7, 2).Value = "Jio Financial Services Limited"
$ws.Cells.Item(7, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(7, 4).Value = "Reducing Consistently"
$ws.Cells.Item(7, 5).Value = 6.436456
$ws.Cells.Item(7, 6).Value = 7.058397
$ws.Cells.Item(7, 7).Value = 7.317609
$ws.Cells.Item(7, 8).Value = -0.6219410000000005
$ws.Cells.Item(7, 9).Value = -0.8811530000000003

# Row 8: INE795G01014 - HDFC Life Insurance Co Ltd
$ws.Cells.Item(8, 1).Value = "INE795G01014"
$ws.Cells.Item(8, 2).Value = "HDFC Life Insurance Co Ltd"
$ws.Cells.Item(8, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(8, 4).Value = "Adding Consistently"
$ws.Cells.Item(8, 5).Value = 5.025197
$ws.Cells.Item(8, 6).Value = 4.877629
$ws.Cells.Item(8, 7).Value = 2.369493
$ws.Cells.Item(8, 8).Value = 0.1475680000000006
$ws.Cells.Item(8, 9).Value = 2.655704000000001

# Row 9: INE216A01030 - Britannia Industries Limited
$ws.Cells.Item(9, 1).Value = "INE216A01030"
$ws.Cells.Item(9, 2).Value = "Britannia Industries Limited"
$ws.Cells.Item(9, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(9, 4).Value = "Adding Consistently"
$ws.Cells.Item(9, 5).Value = 4.949731
$ws.Cells.Item(9, 6).Value = 4.819866
$ws.Cells.Item(9, 7).Value = 4.30574
$ws.Cells.Item(9, 8).Value = 0.1298649999999997
$ws.Cells.Item(9, 9).Value = 0.6439909999999998

# Row 10: INE0J1Y01017 - Life Insurance Corporation Of India
$ws.Cells.Item(10, 1).Value = "INE0J1Y01017"
$ws.Cells.Item(10, 2).Value = "Life Insurance Corporation Of India"
$ws.Cells.Item(10, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(10, 4).Value = "Adding"
$ws.Cells.Item(10, 5).Value = 4.824868
$ws.Cells.Item(10, 6).Value = 4.733788
$ws.Cells.Item(10, 7).Value = 4.937727
$ws.Cells.Item(10, 8).Value = 0.09108000000000072
$ws.Cells.Item(10, 9).Value = -0.1128589999999994

# Row 11: INE245A01021 - Tata Power Company Limited
$ws.Cells.Item(11, 1).Value = "INE245A01021"
$ws.Cells.Item(11, 2).Value = "Tata Power Company Limited"
$ws.Cells.Item(11, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(11, 4).Value = "Adding"
$ws.Cells.Item(11, 5).Value = 4.596682
$ws.Cells.Item(11, 6).Value = 4.507465
$ws.Cells.Item(11, 7).Value = 4.791926
$ws.Cells.Item(11, 8).Value = 0.08921700000000055
$ws.Cells.Item(11, 9).Value = -0.1952439999999998

# Row 12: INE090A01021 - ICICI Bank Limited
$ws.Cells.Item(12, 1).Value = "INE090A01021"
$ws.Cells.Item(12, 2).Value = "ICICI Bank Limited"
$ws.Cells.Item(12, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(12, 4).Value = "Adding Consistently"
$ws.Cells.Item(12, 5).Value = 4.547104
$ws.Cells.Item(12, 6).Value = 0.953632
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(12, 8).Value = 3.593472
$ws.Cells.Item(12, 9).Value = 4.547104

# Row 13: INE917I01010 - Bajaj Auto Limited
$ws.Cells.Item(13, 1).Value = "INE917I01010"
$ws.Cells.Item(13, 2).Value = "Bajaj Auto Limited"
$ws.Cells.Item(13, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(13, 4).Value = "Adding Consistently"
$ws.Cells.Item(13, 5).Value = 4.218914
$ws.Cells.Item(13, 6).Value = 2.04832
$ws.Cells.Item(13, 7).Value = 1.943084
$ws.Cells.Item(13, 8).Value = 2.170594
$ws.Cells.Item(13, 9).Value = 2.27583

# Row 14: INE047A01021 - Grasim Industries Ltd
$ws.Cells.Item(14, 1).Value = "INE047A01021"
$ws.Cells.Item(14, 2).Value = "Grasim Industries Ltd"
$ws.Cells.Item(14, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(14, 4).Value = "Adding Consistently"
$ws.Cells.Item(14, 5).Value = 3.598365
$ws.Cells.Item(14, 6).Value = 3.416975
$ws.Cells.Item(14, 7).Value = 3.481114
$ws.Cells.Item(14, 8).Value = 0.1813899999999999
$ws.Cells.Item(14, 9).Value = 0.117251

# Row 15: INE364U01010 - Adani Green Energy Limited
$ws.Cells.Item(15, 1).Value = "INE364U01010"
$ws.Cells.Item(15, 2).Value = "Adani Green Energy Limited"
$ws.Cells.Item(15, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(15, 4).Value = "Reducing Consistently"
$ws.Cells.Item(15, 5).Value = 3.351149
$ws.Cells.Item(15, 6).Value = 3.776668
$ws.Cells.Item(15, 7).Value = 4.227279
$ws.Cells.Item(15, 8).Value = -0.425519
$ws.Cells.Item(15, 9).Value = -0.8761300000000003

# Row 16: INE271C01023 - DLF Limited
$ws.Cells.Item(16, 1).Value = "INE271C01023"
$ws.Cells.Item(16, 2).Value = "DLF Limited"
$ws.Cells.Item(16, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(16, 4).Value = "Reducing Consistently"
$ws.Cells.Item(16, 5).Value = 2.55428
$ws.Cells.Item(16, 6).Value = 2.613307
$ws.Cells.Item(16, 7).Value = 2.865513
$ws.Cells.Item(16, 8).Value = -0.05902699999999994
$ws.Cells.Item(16, 9).Value = -0.3112330000000001

# Row 17: INE423A01024 - Adani Enterprises Limited
$ws.Cells.Item(17, 1).Value = "INE423A01024"
$ws.Cells.Item(17, 2).Value = "Adani Enterprises Limited"
$ws.Cells.Item(17, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(17, 4).Value = "Reducing Consistently"
$ws.Cells.Item(17, 5).Value = 1.976361
$ws.Cells.Item(17, 6).Value = 2.073086
$ws.Cells.Item(17, 7).Value = 2.288813
$ws.Cells.Item(17, 8).Value = -0.09672499999999995
$ws.Cells.Item(17, 9).Value = -0.3124520000000002

# Row 18: INE180C01042 - Capri Global Capital Limited
$ws.Cells.Item(18, 1).Value = "INE180C01042"
$ws.Cells.Item(18, 2).Value = "Capri Global Capital Limited"
$ws.Cells.Item(18, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(18, 4).Value = "Adding"
$ws.Cells.Item(18, 5).Value = 1.891517
$ws.Cells.Item(18, 6).Value = 1.862128
$ws.Cells.Item(18, 7).Value = 2.060889
$ws.Cells.Item(18, 8).Value = 0.02938899999999989
$ws.Cells.Item(18, 9).Value = -0.1693720000000001

# Row 19: INE237A01036 - Kotak Mahindra Bank Limited
$ws.Cells.Item(19, 1).Value = "INE237A01036"
$ws.Cells.Item(19, 2).Value = "Kotak Mahindra Bank Limited"
$ws.Cells.Item(19, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(19, 4).Value = "Fresh Entry"
$ws.Cells.Item(19, 5).Value = 1.815095
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(19, 8).Value = 1.815095
$ws.Cells.Item(19, 9).Value = 1.815095

# Row 20: INE019C01026 - Himadri Speciality Chemical Limited
$ws.Cells.Item(20, 1).Value = "INE019C01026"
$ws.Cells.Item(20, 2).Value = "Himadri Speciality Chemical Limited"
$ws.Cells.Item(20, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(20, 4).Value = "Reducing"
$ws.Cells.Item(20, 5).Value = 0.927695
$ws.Cells.Item(20, 6).Value = 0.929448
$ws.Cells.Item(20, 7).Value = 0.915729
$ws.Cells.Item(20, 8).Value = -0.001753000000000005
$ws.Cells.Item(20, 9).Value = 0.01196600000000003

# Row 21: INE781S01027 - Ventive Hospitality Limited
$ws.Cells.Item(21, 1).Value = "INE781S01027"
$ws.Cells.Item(21, 2).Value = "Ventive Hospitality Limited"
$ws.Cells.Item(21, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(21, 4).Value = "Adding Consistently"
$ws.Cells.Item(21, 5).Value = 0.882671
$ws.Cells.Item(21, 6).Value = 0.854077
$ws.Cells.Item(21, 7).Value = 0.824463
$ws.Cells.Item(21, 8).Value = 0.02859400000000001
$ws.Cells.Item(21, 9).Value = 0.05820800000000004

# Row 22: INE040A01034 - HDFC Bank Limited
$ws.Cells.Item(22, 1).Value = "INE040A01034"
$ws.Cells.Item(22, 2).Value = "HDFC Bank Limited"
$ws.Cells.Item(22, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(22, 4).Value = "Fresh Entry"
$ws.Cells.Item(22, 5).Value = 0.748248
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(22, 7).Value = 0
$ws.Cells.Item(22, 8).Value = 0.748248
$ws.Cells.Item(22, 9).Value = 0.748248

# Row 23: INE075A01022 - Wipro Ltd
$ws.Cells.Item(23, 1).Value = "INE075A01022"
$ws.Cells.Item(23, 2).Value = "Wipro Ltd"
$ws.Cells.Item(23, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(23, 4).Value = "Reducing"
$ws.Cells.Item(23, 5).Value = 0.601453
$ws.Cells.Item(23, 6).Value = 0.632489
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = -0.03103599999999995
$ws.Cells.Item(23, 9).Value = 0.601453

# Row 24: INE522F01014 - Coal India Ltd
$ws.Cells.Item(24, 1).Value = "INE522F01014"
$ws.Cells.Item(24, 2).Value = "Coal India Ltd"
$ws.Cells.Item(24, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(24, 4).Value = "Fresh Entry"
$ws.Cells.Item(24, 5).Value = 0.552926
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(24, 7).Value = 0
$ws.Cells.Item(24, 8).Value = 0.552926
$ws.Cells.Item(24, 9).Value = 0.552926

# Row 25: INE281B01032 - Lloyds Metals And Energy Limited
$ws.Cells.Item(25, 1).Value = "INE281B01032"
$ws.Cells.Item(25, 2).Value = "Lloyds Metals And Energy Limited"
$ws.Cells.Item(25, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(25, 4).Value = "Reducing Consistently"
$ws.Cells.Item(25, 5).Value = 0.400783
$ws.Cells.Item(25, 6).Value = 0.453962
$ws.Cells.Item(25, 7).Value = 0.446805
$ws.Cells.Item(25, 8).Value = -0.05317899999999998
$ws.Cells.Item(25, 9).Value = -0.04602200000000001

# Row 26: INE101I01011 - Afcons Infrastructure Limited
$ws.Cells.Item(26, 1).Value = "INE101I01011"
$ws.Cells.Item(26, 2).Value = "Afcons Infrastructure Limited"
$ws.Cells.Item(26, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(26, 4).Value = "Reducing Consistently"
$ws.Cells.Item(26, 5).Value = 0.346437
$ws.Cells.Item(26, 6).Value = 0.382479
$ws.Cells.Item(26, 7).Value = 0.44146
$ws.Cells.Item(26, 8).Value = -0.03604200000000002
$ws.Cells.Item(26, 9).Value = -0.09502300000000002

# Row 27: INE069I01010 - Embassy Developments Limited
$ws.Cells.Item(27, 1).Value = "INE069I01010"
$ws.Cells.Item(27, 2).Value = "Embassy Developments Limited"
$ws.Cells.Item(27, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(27, 4).Value = "Adding"
$ws.Cells.Item(27, 5).Value = 0.292822
$ws.Cells.Item(27, 6).Value = 0.246078
$ws.Cells.Item(27, 7).Value = 0.387239
$ws.Cells.Item(27, 8).Value = 0.04674400000000004
$ws.Cells.Item(27, 9).Value = -0.09441699999999997

# Row 28: INE0QN801017 - Krystal Integrated Services Limited
$ws.Cells.Item(28, 1).Value = "INE0QN801017"
$ws.Cells.Item(28, 2).Value = "Krystal Integrated Services Limited"
$ws.Cells.Item(28, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(28, 4).Value = "Adding Consistently"
$ws.Cells.Item(28, 5).Value = 0.269236
$ws.Cells.Item(28, 6).Value = 0.219008
$ws.Cells.Item(28, 7).Value = 0.254633
$ws.Cells.Item(28, 8).Value = 0.05022799999999997
$ws.Cells.Item(28, 9).Value = 0.01460299999999998

# Row 29: INE696F01016 - Juniper Hotels Limited
$ws.Cells.Item(29, 1).Value = "INE696F01016"
$ws.Cells.Item(29, 2).Value = "Juniper Hotels Limited"
$ws.Cells.Item(29, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(29, 4).Value = "Reducing Consistently"
$ws.Cells.Item(29, 5).Value = 0.258687
$ws.Cells.Item(29, 6).Value = 0.274051
$ws.Cells.Item(29, 7).Value = 0.295133
$ws.Cells.Item(29, 8).Value = -0.01536399999999999
$ws.Cells.Item(29, 9).Value = -0.03644599999999998

# Row 30: INE192B01031 - Welspun Living Limited
$ws.Cells.Item(30, 1).Value = "INE192B01031"
$ws.Cells.Item(30, 2).Value = "Welspun Living Limited"
$ws.Cells.Item(30, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(30, 4).Value = "Adding Consistently"
$ws.Cells.Item(30, 5).Value = 0.155441
$ws.Cells.Item(30, 6).Value = 0.154964
$ws.Cells.Item(30, 7).Value = 0.153362
$ws.Cells.Item(30, 8).Value = 0.0004770000000000052
$ws.Cells.Item(30, 9).Value = 0.002078999999999998

# Row 31: INE776C01039 - GMR Airports Limited
$ws.Cells.Item(31, 1).Value = "INE776C01039"
$ws.Cells.Item(31, 2).Value = "GMR Airports Limited"
$ws.Cells.Item(31, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(31, 4).Value = "Complete Exit"
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(31, 6).Value = 1.46356
$ws.Cells.Item(31, 7).Value = 1.31335
$ws.Cells.Item(31, 8).Value = -1.46356
$ws.Cells.Item(31, 9).Value = -1.31335

# Row 32: INE669C01036 - Tech Mahindra Limited
$ws.Cells.Item(32, 1).Value = "INE669C01036"
$ws.Cells.Item(32, 2).Value = "Tech Mahindra Limited"
$ws.Cells.Item(32, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(32, 4).Value = "Complete Exit"
$ws.Cells.Item(32, 5).Value = 0
$ws.Cells.Item(32, 6).Value = 3.033516
$ws.Cells.Item(32, 7).Value = 0
$ws.Cells.Item(32, 8).Value = -3.033516
$ws.Cells.Item(32, 9).Value = 0

# Row 33: INE029A01011 - Bharat Petroleum Corp Ltd
$ws.Cells.Item(33, 1).Value = "INE029A01011"
$ws.Cells.Item(33, 2).Value = "Bharat Petroleum Corp Ltd"
$ws.Cells.Item(33, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(33, 4).Value = "Complete Exit"
$ws.Cells.Item(33, 5).Value = 0
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(33, 7).Value = 1.285099
$ws.Cells.Item(33, 8).Value = 0
$ws.Cells.Item(33, 9).Value = -1.285099

# Row 34: INE044A01036 - Sun Pharmaceutical Industries Limited
$ws.Cells.Item(34, 1).Value = "INE044A01036"
$ws.Cells.Item(34, 2).Value = "Sun Pharmaceutical Industries Limited"
$ws.Cells.Item(34, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(34, 4).Value = "Complete Exit"
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(34, 7).Value = 0.601787
$ws.Cells.Item(34, 8).Value = 0
$ws.Cells.Item(34, 9).Value = -0.601787

# Row 35: INE860A01027 - HCL Technologies Limited
$ws.Cells.Item(35, 1).Value = "INE860A01027"
$ws.Cells.Item(35, 2).Value = "HCL Technologies Limited"
$ws.Cells.Item(35, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(35, 4).Value = "Complete Exit"
$ws.Cells.Item(35, 5).Value = 0
$ws.Cells.Item(35, 6).Value = 0.263262
$ws.Cells.Item(35, 7).Value = 0
$ws.Cells.Item(35, 8).Value = -0.263262
$ws.Cells.Item(35, 9).Value = 0

# Row 36: INE296A01032 - Bajaj Finance Limited
$ws.Cells.Item(36, 1).Value = "INE296A01032"
$ws.Cells.Item(36, 2).Value = "Bajaj Finance Limited"
$ws.Cells.Item(36, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(36, 4).Value = "Complete Exit"
$ws.Cells.Item(36, 5).Value = 0
$ws.Cells.Item(36, 6).Value = 1.298806
$ws.Cells.Item(36, 7).Value = 1.367956
$ws.Cells.Item(36, 8).Value = -1.298806
$ws.Cells.Item(36, 9).Value = -1.367956

# Row 37: INE326A01037 - Lupin Limited
$ws.Cells.Item(37, 1).Value = "INE326A01037"
$ws.Cells.Item(37, 2).Value = "Lupin Limited"
$ws.Cells.Item(37, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(37, 4).Value = "Complete Exit"
$ws.Cells.Item(37, 5).Value = 0
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(37, 7).Value = 1.341811
$ws.Cells.Item(37, 8).Value = 0
$ws.Cells.Item(37, 9).Value = -1.341811

# Row 38: INE285K01026 - Techno Electric & Engineering Co Ltd
$ws.Cells.Item(38, 1).Value = "INE285K01026"
$ws.Cells.Item(38, 2).Value = "Techno Electric & Engineering Co Ltd"
$ws.Cells.Item(38, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(38, 4).Value = "Complete Exit"
$ws.Cells.Item(38, 5).Value = 0
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(38, 7).Value = 0.047635
$ws.Cells.Item(38, 8).Value = 0
$ws.Cells.Item(38, 9).Value = -0.047635

# Row 39: INE059A01026 - Cipla Limited
$ws.Cells.Item(39, 1).Value = "INE059A01026"
$ws.Cells.Item(39, 2).Value = "Cipla Limited"
$ws.Cells.Item(39, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(39, 4).Value = "Complete Exit"
$ws.Cells.Item(39, 5).Value = 0
$ws.Cells.Item(39, 6).Value = 0.2989
$ws.Cells.Item(39, 7).Value = 0.295937
$ws.Cells.Item(39, 8).Value = -0.2989
$ws.Cells.Item(39, 9).Value = -0.295937

# Row 40: INE009A01021 - Infosys Limited
$ws.Cells.Item(40, 1).Value = "INE009A01021"
$ws.Cells.Item(40, 2).Value = "Infosys Limited"
$ws.Cells.Item(40, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(40, 4).Value = "Complete Exit"
$ws.Cells.Item(40, 5).Value = 0
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(40, 7).Value = 1.388588
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 9).Value = -1.388588

# Row 41: INE062A01020 - State Bank of India
$ws.Cells.Item(41, 1).Value = "INE062A01020"
$ws.Cells.Item(41, 2).Value = "State Bank of India"
$ws.Cells.Item(41, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(41, 4).Value = "Complete Exit"
$ws.Cells.Item(41, 5).Value = 0
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(41, 7).Value = 3.810277
$ws.Cells.Item(41, 8).Value = 0
$ws.Cells.Item(41, 9).Value = -3.810277

# Row 42: INE192A01025 - Tata Consumer Products Ltd
$ws.Cells.Item(42, 1).Value = "INE192A01025"
$ws.Cells.Item(42, 2).Value = "Tata Consumer Products Ltd"
$ws.Cells.Item(42, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(42, 4).Value = "Complete Exit"
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(42, 6).Value = 0.981754
$ws.Cells.Item(42, 7).Value = 0
$ws.Cells.Item(42, 8).Value = -0.981754
$ws.Cells.Item(42, 9).Value = 0

# Row 43: INE134E01011 - Power Finance Corporation Ltd.
$ws.Cells.Item(43, 1).Value = "INE134E01011"
$ws.Cells.Item(43, 2).Value = "Power Finance Corporation Ltd."
$ws.Cells.Item(43, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(43, 4).Value = "Complete Exit"
$ws.Cells.Item(43, 5).Value = 0
$ws.Cells.Item(43, 6).Value = 0
$ws.Cells.Item(43, 7).Value = 1.359737
$ws.Cells.Item(43, 8).Value = 0
$ws.Cells.Item(43, 9).Value = -1.359737

# Row 44: INE129A01019 - GAIL (India) Limited
$ws.Cells.Item(44, 1).Value = "INE129A01019"
$ws.Cells.Item(44, 2).Value = "GAIL (India) Limited"
$ws.Cells.Item(44, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(44, 4).Value = "Complete Exit"
$ws.Cells.Item(44, 5).Value = 0
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(44, 7).Value = 1.026296
$ws.Cells.Item(44, 8).Value = 0
$ws.Cells.Item(44, 9).Value = -1.026296

# Row 45: INE237A01028 - Kotak Mahindra Bank Limited
$ws.Cells.Item(45, 1).Value = "INE237A01028"
$ws.Cells.Item(45, 2).Value = "Kotak Mahindra Bank Limited"
$ws.Cells.Item(45, 3).Value = "quant ELSS Tax Saver Fund"
$ws.Cells.Item(45, 4).Value = "Complete Exit"
$ws.Cells.Item(45, 5).Value = 0
$ws.Cells.Item(45, 6).Value = 1.853137
$ws.Cells.Item(45, 7).Value = 0.775955
$ws.Cells.Item(45, 8).Value = -1.853137
$ws.Cells.Item(45, 9).Value = -0.775955
